$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 1: "No preprocessing except standardisation" (rows 8-10) ---
# Row 8 (Gamma=1): value changes only, styles unchanged
$ws.Range("B8").Value = 0.5
$ws.Range("C8").Value = 0.5
$ws.Range("D8").Value = 0.5

# Row 10 (Gamma=3): value changes; C10/D10 also gain the bold style used by B10
$ws.Range("B10").Value = 0.74875999999999998
$ws.Range("C10").Value = 0.72897999999999996
$ws.Range("D10").Value = 0.66620000000000001
$ws.Range("C10").Font.Bold = $true
$ws.Range("D10").Font.Bold = $true

# --- Section 2: "Resolving -999s" (rows 15-17) ---
# Row 15 (Gamma=1): value changes only
$ws.Range("B15").Value = 0.5
$ws.Range("C15").Value = 0.5
$ws.Range("D15").Value = 0.5

# Row 17 (Gamma=3): value changes; C17/D17 also gain the bold style used by B17
$ws.Range("B17").Value = 0.75004000000000004
$ws.Range("C17").Value = 0.72696000000000005
$ws.Range("D17").Value = 0.66679999999999995
$ws.Range("C17").Font.Bold = $true
$ws.Range("D17").Font.Bold = $true

# --- Section 3: "Resolving -999s and removing outliers" (rows 22-24) ---
$ws.Range("B22").Value = 0.5
$ws.Range("B23").Value = 0
$ws.Range("B24").Value = 0.74875999999999998

# --- Section 4: rows 31, 33 ---
$ws.Range("B31").Value = 0.5
$ws.Range("B33").Value = 0.73834

# --- Section 5: rows 38, 40 ---
$ws.Range("B38").Value = 0.5
$ws.Range("B40").Value = 0.73877999999999999

# --- Section 6: rows 45-47 ---
$ws.Range("B45").Value = 0
$ws.Range("B46").Value = 0.5
$ws.Range("B47").Value = 0.72396000000000005

# --- Section 7: rows 51-53 ---
$ws.Range("B51").Value = 0
$ws.Range("B52").Value = 0.5
$ws.Range("B53").Value = 0.72262000000000004

# --- Sheet view state: scroll position and selection ---
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D58").Select()
